$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the license-related values in row 2
$ws.Range("J2").Value = "katronke105991"
$ws.Range("K2").Value = "Killadi12301591"

# Widen column J (10th column) to match new width (~20.54 chars; the
# engine quantizes ColumnWidth to a 1/6-character pixel grid, so
# 118/6 is the closest achievable value to the recorded 20.54296875)
$ws.Columns.Item(10).ColumnWidth = 19.666666666666668

# Move the active selection to K11, as recorded when the sheet was saved
$ws.Range("K11").Select()
